$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data per latest scrape.
# Column D (Price) values that look like plain numbers must be forced to
# text so Excel keeps the original formatted-string representation
# (e.g. preserves trailing zeros / thousand separators) instead of
# converting them to numeric cells.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "64.803.32"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -5.63%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.354.92"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -7.73%  "
$ws.Range("E4").Value = "  +0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "182.59"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -9.77%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "526.97"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -9.69%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.596"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -4.63%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.347.98"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -7.77%  "
$ws.Range("E9").Value = "  -0.02%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.615"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -10.89%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "56.63"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -7.30%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.130"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -13.92%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000250"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -13.56%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "9.11"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -10.25%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.891.95"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -7.31%  "
$ws.Range("E16").Value = "  -4.32%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.358.07"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -7.46%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "64.552.85"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -5.68%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "17.21"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -10.48%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "10.90"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -13.30%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.954"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -11.78%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "368.32"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -9.74%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "80.46"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -6.60%  "
$ws.Range("E24").Value = "  -14.13%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "10.65"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -17.96%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "3.69"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -7.81%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "5.85"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -5.00%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.61"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -11.73%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "11.21"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -12.32%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "8.36"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -12.13%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "665.35"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "29.12"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -8.90%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "6.69"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -15.41%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "60.66"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -5.06%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "11.01"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -11.07%  "
$ws.Range("E36").Value = "  -10.73%  "
$ws.Range("E37").Value = "  -0.12%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "36.00"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -14.67%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.373"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -11.54%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -7.34%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.793.71"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -12.78%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.72"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -16.05%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.60"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -9.49%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0613"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -21.66%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0384"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -8.77%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.31"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -15.46%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.125"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -6.03%  "
$ws.Range("E49").Value = "  -1.84%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.80"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -10.03%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.55"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -6.86%  "
